# Generate Report for Handoff
#
# The localization-status report records, per target locale, the most
# recent "Latest Handoff Datetime" for each source file. Two source
# files (9692b87b-...md and e02d26b1-...md) were just handed off again,
# so their "Latest Handoff Datetime" cells (column D, rows 12 and 14)
# need to be refreshed to the new handoff timestamp on each locale sheet.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("D12").Value = "2016-03-09 07:54:13"
$zhcn.Range("D14").Value = "2016-03-09 07:54:13"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("D12").Value = "2016-03-09 07:54:17"
$dede.Range("D14").Value = "2016-03-09 07:54:17"
